$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumValue($range, $value) {
    $range.Value = $value
}

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-NumValue $ws.Range("B2") 44030
Set-NumValue $ws.Range("C2") 27987
Set-NumValue $ws.Range("D2") 749
Set-NumValue $ws.Range("E2") 3222
Set-NumValue $ws.Range("G2") 11.51
Set-NumValue $ws.Range("H2") 12.68

# Row 3
Set-NumValue $ws.Range("B3") 44030
Set-NumValue $ws.Range("C3") 28633
Set-NumValue $ws.Range("D3") 251

# Row 4
Set-NumValue $ws.Range("B4") 44030
Set-TextValue $ws.Range("C4") "217895"
Set-TextValue $ws.Range("D4") "18771"
Set-NumValue $ws.Range("E4") 33686
Set-NumValue $ws.Range("F4") 5251
Set-NumValue $ws.Range("G4") 30.1
Set-NumValue $ws.Range("K4") 111896
Set-NumValue $ws.Range("L4") 17221

# Row 7
Set-NumValue $ws.Range("B7") 44030
Set-NumValue $ws.Range("C7") 76336
Set-NumValue $ws.Range("D7") 838
Set-NumValue $ws.Range("E7") 15111
Set-NumValue $ws.Range("F7") 300
Set-NumValue $ws.Range("G7") 19.8
Set-NumValue $ws.Range("H7") 35.8

# Row 8
Set-NumValue $ws.Range("B8") 44030
Set-TextValue $ws.Range("C8") "33332"
Set-TextValue $ws.Range("D8") "243"
Set-TextValue $ws.Range("E8") "830"
Set-TextValue $ws.Range("F8") "3"
Set-NumValue $ws.Range("H8") 1.23

# Row 9
Set-NumValue $ws.Range("B9") 44030
Set-NumValue $ws.Range("C9") 22184
Set-NumValue $ws.Range("D9") 667
Set-NumValue $ws.Range("E9") 2343
Set-NumValue $ws.Range("F9") 26
Set-NumValue $ws.Range("G9") 15.49
Set-NumValue $ws.Range("H9") 4.2
Set-NumValue $ws.Range("K9") 15127
Set-NumValue $ws.Range("L9") 622

# Row 10
Set-NumValue $ws.Range("C10") 32533
Set-NumValue $ws.Range("D10") 357
Set-NumValue $ws.Range("E10") 6918
Set-NumValue $ws.Range("F10") 93
Set-NumValue $ws.Range("H10") 26.2
Set-NumValue $ws.Range("K10") 28197
Set-NumValue $ws.Range("L10") 355

# Row 11
Set-NumValue $ws.Range("B11") 44030
Set-NumValue $ws.Range("C11") 23114
Set-NumValue $ws.Range("D11") 478
Set-NumValue $ws.Range("E11") 847
Set-NumValue $ws.Range("G11") 4.68
Set-NumValue $ws.Range("H11") 3.9
Set-NumValue $ws.Range("K11") 18089
Set-NumValue $ws.Range("L11") 462

# Row 13
Set-NumValue $ws.Range("B13") 44030
Set-NumValue $ws.Range("C13") 16736
Set-NumValue $ws.Range("D13") 569
Set-NumValue $ws.Range("E13") 303

# Row 16
Set-NumValue $ws.Range("B16") 44030
Set-NumValue $ws.Range("C16") 64180
Set-NumValue $ws.Range("D16") 1253
Set-NumValue $ws.Range("E16") 19494
Set-NumValue $ws.Range("F16") 536
Set-NumValue $ws.Range("G16") 43.92
Set-NumValue $ws.Range("H16") 44.67
Set-NumValue $ws.Range("K16") 44381
Set-NumValue $ws.Range("L16") 1200

# Row 17
Set-NumValue $ws.Range("B17") 44029
Set-NumValue $ws.Range("C17") 153041
Set-NumValue $ws.Range("D17") 4084
Set-NumValue $ws.Range("E17") 4094
Set-NumValue $ws.Range("F17") 409
Set-NumValue $ws.Range("G17") 4.69
Set-NumValue $ws.Range("H17") 10.76
Set-NumValue $ws.Range("K17") 87304
Set-NumValue $ws.Range("L17") 3801

# Row 19
Set-NumValue $ws.Range("B19") 44029
Set-NumValue $ws.Range("C19") 41846
Set-NumValue $ws.Range("D19") 1346
Set-NumValue $ws.Range("E19") 19138
Set-NumValue $ws.Range("F19") 671
Set-NumValue $ws.Range("G19") 45.73

# Row 22
Set-NumValue $ws.Range("B22") 44030
Set-NumValue $ws.Range("C22") 2471
Set-NumValue $ws.Range("E22") 13
Set-NumValue $ws.Range("G22") 0.53

# Row 26
Set-NumValue $ws.Range("B26") 44030
Set-NumValue $ws.Range("C26") 39788
Set-NumValue $ws.Range("D26") 1752
Set-NumValue $ws.Range("E26") 1981
Set-NumValue $ws.Range("G26") 6.25
Set-NumValue $ws.Range("H26") 6.96
Set-NumValue $ws.Range("K26") 31715
Set-NumValue $ws.Range("L26") 1695

# Row 27
Set-NumValue $ws.Range("B27") 44030
Set-NumValue $ws.Range("C27") 22481

# Row 28
Set-NumValue $ws.Range("B28") 44030
Set-NumValue $ws.Range("C28") 73098
Set-NumValue $ws.Range("D28") 6039
Set-NumValue $ws.Range("E28") 21215
Set-NumValue $ws.Range("F28") 2406
Set-NumValue $ws.Range("G28") 29.02
Set-NumValue $ws.Range("H28") 39.84

# Row 29
Set-NumValue $ws.Range("B29") 44029
Set-NumValue $ws.Range("C29") 375363
Set-NumValue $ws.Range("D29") 7595
Set-NumValue $ws.Range("E29") 10432
Set-NumValue $ws.Range("F29") 641
Set-NumValue $ws.Range("G29") 4.32
Set-NumValue $ws.Range("H29") 8.69
Set-NumValue $ws.Range("K29") 241390
Set-NumValue $ws.Range("L29") 7376

# Row 30
Set-NumValue $ws.Range("B30") 44030
Set-NumValue $ws.Range("C30") 55654
Set-NumValue $ws.Range("D30") 2627
Set-NumValue $ws.Range("E30") 6493
Set-NumValue $ws.Range("G30") 11.67
Set-NumValue $ws.Range("H30") 14.2

# Row 31
Set-NumValue $ws.Range("B31") 44030
Set-NumValue $ws.Range("C31") 1795
Set-NumValue $ws.Range("D31") 18
Set-NumValue $ws.Range("E31") 42
Set-NumValue $ws.Range("G31") 1.4
Set-NumValue $ws.Range("K31") 3010
Set-NumValue $ws.Range("L31") 36

# Row 32
Set-NumValue $ws.Range("B32") 44030
Set-NumValue $ws.Range("C32") 41485
Set-NumValue $ws.Range("D32") 843
Set-NumValue $ws.Range("E32") 6721
Set-NumValue $ws.Range("G32") 17.93
Set-NumValue $ws.Range("H32") 23.64
Set-NumValue $ws.Range("K32") 37485
Set-NumValue $ws.Range("L32") 829

# Row 33
Set-NumValue $ws.Range("B33") 44030
Set-NumValue $ws.Range("C33") 139872
Set-NumValue $ws.Range("D33") 3168
Set-NumValue $ws.Range("E33") 36504
Set-NumValue $ws.Range("F33") 1470
Set-NumValue $ws.Range("G33") 26.1
Set-NumValue $ws.Range("H33") 46.4

# Row 34
Set-NumValue $ws.Range("B34") 44030
Set-NumValue $ws.Range("C34") 46026
Set-NumValue $ws.Range("D34") 1444
Set-NumValue $ws.Range("E34") 1690
Set-NumValue $ws.Range("G34") 5.42
Set-NumValue $ws.Range("H34") 3.46
Set-NumValue $ws.Range("K34") 31176
Set-NumValue $ws.Range("L34") 1386

# Row 38
Set-NumValue $ws.Range("C38") 38197
Set-NumValue $ws.Range("D38") 789
Set-NumValue $ws.Range("E38") 3188
Set-NumValue $ws.Range("G38") 8.35
Set-NumValue $ws.Range("H38") 4.82

# Row 39
Set-NumValue $ws.Range("B39") 44030
Set-NumValue $ws.Range("C39") 97958
Set-NumValue $ws.Range("D39") 1629
Set-NumValue $ws.Range("E39") 15936
Set-NumValue $ws.Range("F39") 516
Set-NumValue $ws.Range("G39") 23.89
Set-NumValue $ws.Range("H39") 32.8
Set-NumValue $ws.Range("K39") 66695
Set-NumValue $ws.Range("L39") 1573

# Row 40
Set-NumValue $ws.Range("B40") 44030
Set-NumValue $ws.Range("C40") 160610
Set-NumValue $ws.Range("D40") 7290
Set-NumValue $ws.Range("E40") 27009
Set-NumValue $ws.Range("F40") 2011
Set-NumValue $ws.Range("G40") 16.82

# Row 41
Set-NumValue $ws.Range("B41") 44030
Set-NumValue $ws.Range("C41") 14302
Set-NumValue $ws.Range("E41") 177
Set-NumValue $ws.Range("G41") 1.24

# Row 42
Set-NumValue $ws.Range("B42") 44030
Set-NumValue $ws.Range("C42") 45470
Set-NumValue $ws.Range("D42") 1538
Set-NumValue $ws.Range("E42") 9200
Set-NumValue $ws.Range("F42") 151
Set-NumValue $ws.Range("G42") 20.23
Set-NumValue $ws.Range("H42") 9.82

# Row 43
Set-NumValue $ws.Range("B43") 44030
Set-NumValue $ws.Range("C43") 113238
Set-NumValue $ws.Range("D43") 8419
Set-NumValue $ws.Range("E43") 10642
Set-NumValue $ws.Range("F43") 691

# Row 45
Set-NumValue $ws.Range("B45") 44030
Set-NumValue $ws.Range("C45") 32246
Set-NumValue $ws.Range("D45") 1130
Set-NumValue $ws.Range("E45") 7584
Set-NumValue $ws.Range("F45") 380
Set-NumValue $ws.Range("G45") 32.77
Set-NumValue $ws.Range("H45") 36.09
Set-NumValue $ws.Range("K45") 23144
Set-NumValue $ws.Range("L45") 1053
